$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Day 11 (column O) attendance: mark "p" (presente) for every student
# except the two rows that were left blank (rows 8 and 13).
$rows = 3..34 | Where-Object { $_ -ne 8 -and $_ -ne 13 }
foreach ($r in $rows) {
    $ws.Range("O$r").Value = "p"
}

# The "Correos" column (D) is unhidden again, back to the same width as
# the other normal columns (A/B).
$ws.Columns.Item(4).Hidden = $false
$ws.Columns.Item(4).ColumnWidth = 8.29

# Restore the view to the top of the sheet and leave the cursor on O6.
$ws.Range("O6").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 3
